# Add two new rows (49 and 50) of match-stat data to the bottom of the
# "ataque" sheet (Plan1 / ActiveSheet), then move the view/selection to
# match where the author ended up after the edit (bottom of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, columns A..U (21 columns each)
$row49 = @(48, 38, 62, 11, 15, 5, 2, 6, 9, 5, 6, 336, 537, 261, 470, 2, 4, 60, 28, 13, 0)
$row50 = @(49, 61, 39, 27, 17, 7, 5, 13, 12, 14, 5, 502, 322, 448, 262, 10, 3, 30, 18, 29, 8)

$data = New-Object 'object[,]' 2,21
for ($c = 0; $c -lt 21; $c++) {
    $data[0, $c] = $row49[$c]
    $data[1, $c] = $row50[$c]
}

$ws.Range("A49:U50").Value = $data

# Move the view the same way the author's workbook ended up: scrolled down
# so row 31 is at the top and column M is the left-most visible column,
# with the final selection on V50 (one cell past the last used column).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 13
$ws.Range("V50").Select()
